# Apply updated transition-probability values to Sheet1 (Stetson_B team-specific matrix)
# This corresponds to "added team spec time commit pt2": refreshed probability
# values across rows 2-4, 6-13, and 15-19 of the matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{Cell="B2"; Value=0.2480916030534351},
    @{Cell="C2"; Value=0.4580152671755725},
    @{Cell="J2"; Value=0.03053435114503817},
    @{Cell="P2"; Value=0.1870229007633588},
    @{Cell="S2"; Value=0.07633587786259542},
    @{Cell="B3"; Value=0.00819672131147541},
    @{Cell="C3"; Value=0.02459016393442623},
    @{Cell="J3"; Value=0.04098360655737705},
    @{Cell="P3"; Value=0.7868852459016393},
    @{Cell="S3"; Value=0.139344262295082},
    @{Cell="J4"; Value=0.07142857142857142},
    @{Cell="P4"; Value=0.7142857142857143},
    @{Cell="S4"; Value=0.2142857142857143},
    @{Cell="B6"; Value=0.05263157894736842},
    @{Cell="D6"; Value=0.008771929824561403},
    @{Cell="F6"; Value=0.07017543859649122},
    @{Cell="J6"; Value=0.2324561403508772},
    @{Cell="O6"; Value=0.0131578947368421},
    @{Cell="Q6"; Value=0.1710526315789474},
    @{Cell="R6"; Value=0.08771929824561403},
    @{Cell="S6"; Value=0.3640350877192983},
    @{Cell="B7"; Value=0.075},
    @{Cell="D7"; Value=0.055},
    @{Cell="E7"; Value=0.005},
    @{Cell="F7"; Value=0.06},
    @{Cell="J7"; Value=0.125},
    @{Cell="O7"; Value=0.01},
    @{Cell="Q7"; Value=0.16},
    @{Cell="R7"; Value=0.075},
    @{Cell="S7"; Value=0.435},
    @{Cell="B8"; Value=0.06326530612244897},
    @{Cell="D8"; Value=0.01224489795918367},
    @{Cell="F8"; Value=0.05510204081632653},
    @{Cell="J8"; Value=0.1408163265306122},
    @{Cell="O8"; Value=0.02040816326530612},
    @{Cell="Q8"; Value=0.1714285714285714},
    @{Cell="R8"; Value=0.1061224489795918},
    @{Cell="S8"; Value=0.4306122448979592},
    @{Cell="B9"; Value=0.07116104868913857},
    @{Cell="D9"; Value=0.0149812734082397},
    @{Cell="F9"; Value=0.0599250936329588},
    @{Cell="J9"; Value=0.1086142322097378},
    @{Cell="O9"; Value=0.01123595505617977},
    @{Cell="Q9"; Value=0.2134831460674157},
    @{Cell="R9"; Value=0.08239700374531835},
    @{Cell="S9"; Value=0.4382022471910113},
    @{Cell="B10"; Value=0.09042553191489362},
    @{Cell="D10"; Value=0.01595744680851064},
    @{Cell="E10"; Value=0.0007598784194528875},
    @{Cell="F10"; Value=0.06610942249240122},
    @{Cell="J10"; Value=0.1238601823708207},
    @{Cell="O10"; Value=0.01595744680851064},
    @{Cell="Q10"; Value=0.2272036474164134},
    @{Cell="R10"; Value=0.0972644376899696},
    @{Cell="S10"; Value=0.3624620060790273},
    @{Cell="G11"; Value=0.1135531135531136},
    @{Cell="J11"; Value=0.08791208791208792},
    @{Cell="K11"; Value=0.1758241758241758},
    @{Cell="L11"; Value=0.608058608058608},
    @{Cell="S11"; Value=0.01465201465201465},
    @{Cell="G12"; Value=0.8011695906432749},
    @{Cell="J12"; Value=0.1228070175438596},
    @{Cell="K12"; Value=0.01169590643274854},
    @{Cell="L12"; Value=0.03508771929824561},
    @{Cell="S12"; Value=0.02923976608187134},
    @{Cell="G13"; Value=0.7346938775510204},
    @{Cell="J13"; Value=0.2244897959183673},
    @{Cell="S13"; Value=0.04081632653061224},
    @{Cell="F15"; Value=0.0392156862745098},
    @{Cell="H15"; Value=0.1843137254901961},
    @{Cell="I15"; Value=0.06274509803921569},
    @{Cell="J15"; Value=0.3647058823529412},
    @{Cell="K15"; Value=0.05490196078431372},
    @{Cell="M15"; Value=0.01176470588235294},
    @{Cell="O15"; Value=0.06274509803921569},
    @{Cell="S15"; Value=0.2196078431372549},
    @{Cell="F16"; Value=0.01176470588235294},
    @{Cell="H16"; Value=0.2176470588235294},
    @{Cell="I16"; Value=0.05882352941176471},
    @{Cell="J16"; Value=0.4176470588235294},
    @{Cell="K16"; Value=0.1058823529411765},
    @{Cell="M16"; Value=0.01764705882352941},
    @{Cell="O16"; Value=0.05882352941176471},
    @{Cell="S16"; Value=0.1117647058823529},
    @{Cell="F17"; Value=0.01972386587771203},
    @{Cell="H17"; Value=0.1479289940828402},
    @{Cell="I17"; Value=0.1124260355029586},
    @{Cell="J17"; Value=0.4161735700197239},
    @{Cell="K17"; Value=0.106508875739645},
    @{Cell="M17"; Value=0.01380670611439842},
    @{Cell="O17"; Value=0.07297830374753451},
    @{Cell="S17"; Value=0.1104536489151874},
    @{Cell="F18"; Value=0.02542372881355932},
    @{Cell="H18"; Value=0.173728813559322},
    @{Cell="I18"; Value=0.1059322033898305},
    @{Cell="J18"; Value=0.3813559322033898},
    @{Cell="K18"; Value=0.09745762711864407},
    @{Cell="M18"; Value=0.01694915254237288},
    @{Cell="O18"; Value=0.1016949152542373},
    @{Cell="S18"; Value=0.09745762711864407},
    @{Cell="F19"; Value=0.01357466063348416},
    @{Cell="H19"; Value=0.220211161387632},
    @{Cell="I19"; Value=0.1206636500754148},
    @{Cell="J19"; Value=0.3423831070889894},
    @{Cell="K19"; Value=0.08521870286576169},
    @{Cell="M19"; Value=0.0248868778280543},
    @{Cell="O19"; Value=0.07239819004524888},
    @{Cell="S19"; Value=0.1206636500754148}

)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host ("Updated {0} cells" -f $updates.Count)
